$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# --- Copy formatting for header row (991) from an existing MONDAY header row (row 9) ---
$ws.Range("A9:F9").Copy()
$ws.Range("A991:F991").PasteSpecial(-4122)
$ws.Cells.Item(991,2).Value = 'MONDAY'

# --- Copy formatting for first two data rows after header (992, 993) from row 408s pattern ---
$ws.Range("A408:F408").Copy()
$ws.Range("A992:F993").PasteSpecial(-4122)

# --- Populate data rows 992-1011 ---
$ws.Range("A992").Value = 'Demo'
$ws.Range("B992").Value = 42744
$ws.Range("C992").Value = '1620'
$ws.Range("D992").Value = 'CLH'
$ws.Range("E992").Value = 'I'
$ws.Range("F992").Value = 'Make sure client is okay.'

$ws.Range("A993").Value = 'AV Shutdown'
$ws.Range("B993").Value = 42744
$ws.Range("C993").Value = '1630'
$ws.Range("D993").Value = 'LAS'
$ws.Range("E993").Value = 'B'
$ws.Range("F993").Value = 'Make sure neck mic goes back to drawer and log off touchscreen.'

$ws.Range("A994").Value = 'Setup Mic'
$ws.Range("B994").Value = 42744
$ws.Range("C994").Value = '1715'
$ws.Range("D994").Value = 'LAS'
$ws.Range("E994").Value = 'C'
$ws.Range("F994").Value = 'Take cart with mixer, 2 wireless mics and 2 mic stands from Lassonde 1011 storeroom (across from Lassonde A). Go to Lassonde C classroom (class starts at 5:30 pm but be there early in case previous class ends early). '

$ws.Range("A995").Value = 'Other'
$ws.Range("B995").Value = 42744
$ws.Range("C995").Value = '1715'
$ws.Range("D995").Value = 'LAS'
$ws.Range("E995").Value = 'C'
$ws.Range("F995").Value = 'Log in as 5065*0 on touchscreen. (First level bar is your wireless handheld mic volume). Plug in mic cable from output of mixer to mic input on podium (XLR jack just above VHS machine in podium). Ramp up volume a bit on "Microphone 2" on touchscreen to medium volume to get level.'

$ws.Range("A996").Value = 'Other'
$ws.Range("B996").Value = 42744
$ws.Range("C996").Value = '1715'
$ws.Range("D996").Value = 'LAS'
$ws.Range("E996").Value = 'C'
$ws.Range("F996").Value = 'Plug in power cord from cart on to power outlet on left side of podium (to left of document camera). Turn on mixer. Turn on wireless microphone receivers on cart (NOTE: DO NOT PRESS "SYNC" BUTTON" - POWER BUTTON IS FIRST BUTTON TO THE RIGHT ON RECEIVER). '

$ws.Range("A997").Value = 'Other'
$ws.Range("B997").Value = 42744
$ws.Range("C997").Value = '1715'
$ws.Range("D997").Value = 'LAS'
$ws.Range("E997").Value = 'C'
$ws.Range("F997").Value = 'Press "MUTE" button on wireless mics to turn on mics.'

$ws.Range("A998").Value = 'Other'
$ws.Range("B998").Value = 42744
$ws.Range("C998").Value = '1715'
$ws.Range("D998").Value = 'LAS'
$ws.Range("E998").Value = 'C'
$ws.Range("F998").Value = 'Once volumes are set, place one mic stand with mic halfway up aisle on right and one mic stand with mic halfway up aisle on left. Demo volume controls to prof. and demo PC. Leave microphone bags with milk carton on cart in room. PLEASE FIND OUT END TIME OF CLASS FROM PROF. AND TELL MASI AS MICROPHONES ARE EXPENSIVE. TELL PROF. TO STAY WITH MICS UNTIL THEY ARE PICKED UP. TELL HIM TO CALL ext 55800   WHEN DONE (use phone in classroom).'

$ws.Range("A999").Value = 'AV Shutdown'
$ws.Range("B999").Value = 42744
$ws.Range("C999").Value = '1820'
$ws.Range("D999").Value = 'CLH'
$ws.Range("E999").Value = 'I'
$ws.Range("F999").Value = 'Make sure neck mic goes back to drawer and log off touchscreen.'

$ws.Range("A1000").Value = 'Demo'
$ws.Range("B1000").Value = 42744
$ws.Range("C1000").Value = '1850'
$ws.Range("D1000").Value = 'CLH'
$ws.Range("E1000").Value = 'B'
$ws.Range("F1000").Value = 'Demo and show prof how to log off and return mic, cables, etc. to drawer. '

$ws.Range("A1001").Value = 'AV Shutdown'
$ws.Range("B1001").Value = 42744
$ws.Range("C1001").Value = '1900'
$ws.Range("D1001").Value = 'CLH'
$ws.Range("E1001").Value = 'L'
$ws.Range("F1001").Value = 'PLEASE MAKE SURE CRESTRON GETS LOGGED OFF. WE ARE HAVING PROBLEMS WITH THIS ROOM WHEN IT DOESN''T GET LOGGED OFF.'

$ws.Range("A1002").Value = 'AV Shutdown'
$ws.Range("B1002").Value = 42744
$ws.Range("C1002").Value = '1630'
$ws.Range("D1002").Value = 'LSB'
$ws.Range("E1002").Value = '101'
$ws.Range("F1002").Value = 'Make sure neck mic goes back to drawer and log off touchscreen.'

$ws.Range("A1003").Value = 'AV Shutdown'
$ws.Range("B1003").Value = 42744
$ws.Range("C1003").Value = '1900'
$ws.Range("D1003").Value = 'LSB'
$ws.Range("E1003").Value = '103'
$ws.Range("F1003").Value = 'Make sure neck mic goes back to drawer and log off touchscreen.'

$ws.Range("A1004").Value = 'AV Shutdown'
$ws.Range("B1004").Value = 42744
$ws.Range("C1004").Value = '1900'
$ws.Range("D1004").Value = 'LSB'
$ws.Range("E1004").Value = '105'
$ws.Range("F1004").Value = 'Make sure neck mic goes back to drawer and log off touchscreen.'

$ws.Range("A1005").Value = 'AV Shutdown'
$ws.Range("B1005").Value = 42744
$ws.Range("C1005").Value = '2000'
$ws.Range("D1005").Value = 'LSB'
$ws.Range("E1005").Value = '106'
$ws.Range("F1005").Value = 'Make sure neck mic goes back to drawer and log off touchscreen.'

$ws.Range("A1006").Value = 'AV Shutdown'
$ws.Range("B1006").Value = 42744
$ws.Range("C1006").Value = '1730'
$ws.Range("D1006").Value = 'LSB'
$ws.Range("E1006").Value = '107'
$ws.Range("F1006").Value = 'Make sure neck mic goes back to drawer and log off touchscreen.'

$ws.Range("A1007").Value = 'Pickup Mic'
$ws.Range("B1007").Value = 42744
$ws.Range("C1007").Value = '1850'
$ws.Range("D1007").Value = 'LAS'
$ws.Range("E1007").Value = 'C'
$ws.Range("F1007").Value = 'Pick up 2 wireless mics on stands with cart. Move all equipment on cart - cart has 2 wireless mic receivers and mixer and mic cables. Pick up 2 mic stands - return all equipment to Lassonde 1011 storeroom (across the hall from Lassonde A). PLEASE PUT 2 WIRELESS MICS IN BAGS PROVIDED IN MILK CARTON ON CART. Very expensive mics - please go early and treat mics with care.'

$ws.Range("A1008").Value = 'Other'
$ws.Range("B1008").Value = 42744
$ws.Range("C1008").Value = '1850'
$ws.Range("D1008").Value = 'LAS'
$ws.Range("E1008").Value = 'C'
$ws.Range("F1008").Value = 'Turn off wireless microphones by pressing "MUTE" button on mics.'

$ws.Range("A1009").Value = 'Other'
$ws.Range("B1009").Value = 42744
$ws.Range("C1009").Value = '1850'
$ws.Range("D1009").Value = 'LAS'
$ws.Range("E1009").Value = 'C'
$ws.Range("F1009").Value = 'Turn off wireless microphone receivers by pressing "POWER" button and not "SYNC" button. '

$ws.Range("A1010").Value = 'Other'
$ws.Range("B1010").Value = 42744
$ws.Range("C1010").Value = '1850'
$ws.Range("D1010").Value = 'LAS'
$ws.Range("E1010").Value = 'C'
$ws.Range("F1010").Value = 'PLEASE BE ON TIME - Prof upset last week when no one came till 7:05 pm and other class was starting.'

$ws.Range("A1011").Value = 'AV Shutdown'
$ws.Range("B1011").Value = 42744
$ws.Range("C1011").Value = '2100'
$ws.Range("D1011").Value = 'LAS'
$ws.Range("E1011").Value = 'A'
$ws.Range("F1011").Value = 'Make sure neck mic goes back to drawer and log off touchscreen.'

# --- Explicit row heights matching wrapped-text auto-fit heights from the source row ---
$ws.Rows.Item(994).RowHeight = 60
$ws.Rows.Item(995).RowHeight = 75
$ws.Rows.Item(996).RowHeight = 75
$ws.Rows.Item(998).RowHeight = 120
$ws.Rows.Item(1000).RowHeight = 30
$ws.Rows.Item(1001).RowHeight = 45
$ws.Rows.Item(1007).RowHeight = 90
$ws.Rows.Item(1008).RowHeight = 30
$ws.Rows.Item(1009).RowHeight = 30
$ws.Rows.Item(1010).RowHeight = 30

# --- Update selection / active cell to match end of log ---
$ws.Range("A1011").Select()

